$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.4458169960283037
$ws.Range("J2").Value = 0.4458169960283037
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 21.18229649518089
$ws.Range("R2").Value = 190.640668456628
$ws.Range("S2").Value = 0.02859981413225822
$ws.Range("T2").Value = 0.02859981413225822

# Row 3
$ws.Range("I3").Value = 0.4458169960283037
$ws.Range("J3").Value = 0.4458169960283037
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("S3").Value = 0.1773951131656029
$ws.Range("T3").Value = 0.1773951131656029

# Row 4
$ws.Range("I4").Value = 0.4458169960283037
$ws.Range("J4").Value = 0.4458169960283037
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 52.61124199329555
$ws.Range("R4").Value = 473.50117793966
$ws.Range("S4").Value = 0.07103440094976643
$ws.Range("T4").Value = 0.07103440094976643

# Row 5
$ws.Range("I5").Value = 0.4458169960283037
$ws.Range("J5").Value = 0.4458169960283037
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 125.0116664089687
$ws.Range("R5").Value = 1125.104997680718
$ws.Range("S5").Value = 0.1687876677806761
$ws.Range("T5").Value = 0.1687876677806761

# Row 6
$ws.Range("G6").Value = 1.399743666666667
$ws.Range("H6").Value = 4.199231
$ws.Range("I6").Value = 0.2598558798146963
$ws.Range("J6").Value = 0.2598558798146962
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 12.34664524073356
$ws.Range("R6").Value = 111.119807166602
$ws.Range("S6").Value = 0.01667013579581636
$ws.Range("T6").Value = 0.01667013579581635

# Row 7
$ws.Range("G7").Value = 1.399743666666667
$ws.Range("H7").Value = 4.199231
$ws.Range("I7").Value = 0.2598558798146963
$ws.Range("J7").Value = 0.2598558798146962
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("Q7").Value = 76.58212460986167
$ws.Range("R7").Value = 689.239121488755
$ws.Range("S7").Value = 0.1033992952649763
$ws.Range("T7").Value = 0.1033992952649763

# Row 8
$ws.Range("G8").Value = 1.399743666666667
$ws.Range("H8").Value = 4.199231
$ws.Range("I8").Value = 0.2598558798146963
$ws.Range("J8").Value = 0.2598558798146962
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 30.66581287413223
$ws.Range("R8").Value = 275.99231586719
$ws.Range("S8").Value = 0.04140422397610781
$ws.Range("T8").Value = 0.0414042239761078

# Row 9
$ws.Range("G9").Value = 1.399743666666667
$ws.Range("H9").Value = 4.199231
$ws.Range("I9").Value = 0.2598558798146963
$ws.Range("J9").Value = 0.2598558798146962
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 72.86625869180966
$ws.Range("R9").Value = 655.796328226287
$ws.Range("S9").Value = 0.09838222477779575
$ws.Range("T9").Value = 0.09838222477779575

# Row 10
$ws.Range("G10").Value = 1.585427
$ws.Range("H10").Value = 4.756281
$ws.Range("I10").Value = 0.294327124157
$ws.Range("J10").Value = 0.294327124157
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 13.98449243974467
$ws.Range("R10").Value = 125.860431957702
$ws.Range("S10").Value = 0.01888151667604407
$ws.Range("T10").Value = 0.01888151667604407

# Row 11
$ws.Range("G11").Value = 1.585427
$ws.Range("H11").Value = 4.756281
$ws.Range("I11").Value = 0.294327124157
$ws.Range("J11").Value = 0.294327124157
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 86.74114480044499
$ws.Range("R11").Value = 780.6703032040049
$ws.Range("S11").Value = 0.1171157536897105
$ws.Range("T11").Value = 0.1171157536897105

# Row 12
$ws.Range("G12").Value = 1.585427
$ws.Range("H12").Value = 4.756281
$ws.Range("I12").Value = 0.294327124157
$ws.Range("J12").Value = 0.294327124157
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 34.73379366907667
$ws.Range("R12").Value = 312.60414302169
$ws.Range("S12").Value = 0.04689671128292443
$ws.Range("T12").Value = 0.04689671128292442

# Row 13
$ws.Range("G13").Value = 1.585427
$ws.Range("H13").Value = 4.756281
$ws.Range("I13").Value = 0.294327124157
$ws.Range("J13").Value = 0.294327124157
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 82.53234979379299
$ws.Range("R13").Value = 742.7911481441369
$ws.Range("S13").Value = 0.111433142508321
$ws.Range("T13").Value = 0.111433142508321
